$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5579
$ws1.Range("F6").Value = 78
$ws1.Range("F8").Value = 920
$ws1.Range("F9").Value = 148
$ws1.Range("F10").Value = 2496
$ws1.Range("F12").Value = 117
$ws1.Range("F14").Value = 75
$ws1.Range("F15").Value = 12
$ws1.Range("F16").Value = 2333
$ws1.Range("F17").Value = 294

# Sheet "全部类型" (sheet4): update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5579
$ws4.Range("F7").Value = 78
$ws4.Range("F10").Value = 920
$ws4.Range("F11").Value = 148
$ws4.Range("F12").Value = 2496
$ws4.Range("F14").Value = 117
$ws4.Range("F17").Value = 75
$ws4.Range("F18").Value = 12
$ws4.Range("F19").Value = 2333
$ws4.Range("F20").Value = 294
